# Renames the three inline logo pictures that live in the document's
# headers/footers:
#   - header1.xml  (first-page header, BTEC logo)   image1.jpg -> image2.jpg
#   - footer1.xml  (first-page footer, Pearson logo) image2.png -> image1.png
#   - footer2.xml  (default footer,    Pearson logo) image2.png -> image1.png
#
# Word only exposes these as InlineShapes hanging off the section's
# Headers/Footers collections -- there is no "rename" dialog command, so we
# drive it the same way a user would: grab the picture and set its .Name.

$d   = $word.ActiveDocument
$sec = $d.Sections.First

function Rename-LogoShape($headerFooter, $newName) {
    if ($headerFooter -eq $null) { return }
    if (-not $headerFooter.Exists) { return }
    $shapes = $headerFooter.Range.InlineShapes
    if ($shapes.Count -lt 1) { return }
    $shapes.Item(1).Name = $newName
}

# First-page header (index 2 of the Headers collection) holds the BTEC logo.
Rename-LogoShape $sec.Headers.Item(2) "image2.jpg"

# Primary/default footer (index 1 of the Footers collection) holds the
# Pearson logo currently named image2.png.
Rename-LogoShape $sec.Footers.Item(1) "image1.png"

# First-page footer (index 2 of the Footers collection) holds the other
# copy of the Pearson logo, also currently named image2.png.
Rename-LogoShape $sec.Footers.Item(2) "image1.png"
